# Junction_Flooding_70.xlsx edit
# - Row 5 values are replaced with a "custom accuracy" (rounded) dataset
# - Row 6 (the old last data row) is removed entirely
# - Sheet dimension shrinks from A1:AH6 to A1:AH5 (handled automatically by
#   the row delete)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 5 with the new, rounded measurement values ---------------
$ws.Range("A5").Value = 40751.4027662037

$ws.Range("B5").Value  = 13.45
$ws.Range("C5").Value  = 9.77
$ws.Range("D5").Value  = 1
$ws.Range("E5").Value  = 29.23
$ws.Range("F5").Value  = 23.6
$ws.Range("G5").Value  = 10.59
$ws.Range("H5").Value  = 44
$ws.Range("I5").Value  = 16.29
$ws.Range("J5").Value  = 7.17
$ws.Range("K5").Value  = 10.5
$ws.Range("L5").Value  = 11.73
$ws.Range("M5").Value  = 12.33
$ws.Range("N5").Value  = 3.38
$ws.Range("O5").Value  = 10.53
$ws.Range("P5").Value  = 14.93
$ws.Range("Q5").Value  = 9
$ws.Range("R5").Value  = 0.78
$ws.Range("S5").Value  = 0.64
$ws.Range("T5").Value  = 152.92
$ws.Range("U5").Value  = 29.53
$ws.Range("V5").Value  = 9.720000000000001
$ws.Range("W5").Value  = 19.71
$ws.Range("X5").Value  = 10.29
$ws.Range("Y5").Value  = 1.72
$ws.Range("Z5").Value  = 20.91
$ws.Range("AA5").Value = 8.58
$ws.Range("AB5").Value = 7.66
$ws.Range("AC5").Value = 9.02
$ws.Range("AD5").Value = 12.31
$ws.Range("AE5").Value = 0.55
$ws.Range("AF5").Value = 40.17
$ws.Range("AG5").Value = 5.4
$ws.Range("AH5").Value = 12.15

# --- Remove the now-obsolete row 6 (shrinks the sheet to A1:AH5) ---------
$ws.Rows.Item(6).Delete()
